# Auto-applied numeric updates to the Leve profit-tracking sheets.
# For each (sheet, cell) pair below, write the new cached value that
# Excel recomputed for that row (currentAveragePrice / LevePrice* / LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$sheetEdits = @{}

$sheetEdits["ALC"] = @(
    @{ Cell = "H15"; Value = 1011.25 },
    @{ Cell = "I15"; Value = 1011.25 },
    @{ Cell = "K15"; Value = 3033.75 },
    @{ Cell = "M15"; Value = -2864.75 },
    @{ Cell = "H86"; Value = 1416.1666 },
    @{ Cell = "I86"; Value = 1124.25 },
    @{ Cell = "K86"; Value = 1124.25 },
    @{ Cell = "M86"; Value = -1.25 },
    @{ Cell = "H89"; Value = 1416.1666 },
    @{ Cell = "I89"; Value = 1124.25 },
    @{ Cell = "K89"; Value = 5621.25 },
    @{ Cell = "M89"; Value = -5.25 },
    @{ Cell = "H100"; Value = 22223664 },
    @{ Cell = "I100"; Value = 28573020 },
    @{ Cell = "J100"; Value = 920 },
    @{ Cell = "K100"; Value = 28573020 },
    @{ Cell = "L100"; Value = 920 },
    @{ Cell = "M100"; Value = -28572479 },
    @{ Cell = "N100"; Value = -2002 },
    @{ Cell = "H106"; Value = 1566 },
    @{ Cell = "I106"; Value = 1566 },
    @{ Cell = "K106"; Value = 1566 },
    @{ Cell = "M106"; Value = -935 },
    @{ Cell = "H115"; Value = 1450.5 },
    @{ Cell = "I115"; Value = 1389.4445 },
    @{ Cell = "K115"; Value = 4168.333500000001 },
    @{ Cell = "M115"; Value = -2601.333500000001 },
    @{ Cell = "H137"; Value = 2072708.2 },
    @{ Cell = "I137"; Value = 2977418 },
    @{ Cell = "J137"; Value = 4800 },
    @{ Cell = "K137"; Value = 8932254 },
    @{ Cell = "L137"; Value = 14400 },
    @{ Cell = "M137"; Value = -8929704 },
    @{ Cell = "N137"; Value = -19500 }
)

$sheetEdits["ARM"] = @(
    @{ Cell = "H61"; Value = 2104.6843 },
    @{ Cell = "I61"; Value = 1326.7778 },
    @{ Cell = "J61"; Value = 2804.8 },
    @{ Cell = "K61"; Value = 1326.7778 },
    @{ Cell = "L61"; Value = 2804.8 },
    @{ Cell = "M61"; Value = -1114.7778 },
    @{ Cell = "N61"; Value = -3228.8 },
    @{ Cell = "H74"; Value = 4702.6294 },
    @{ Cell = "I74"; Value = 6011.9375 },
    @{ Cell = "J74"; Value = 2798.182 },
    @{ Cell = "K74"; Value = 6011.9375 },
    @{ Cell = "L74"; Value = 2798.182 },
    @{ Cell = "M74"; Value = -5137.9375 },
    @{ Cell = "N74"; Value = -4546.182 },
    @{ Cell = "H77"; Value = 4702.6294 },
    @{ Cell = "I77"; Value = 6011.9375 },
    @{ Cell = "J77"; Value = 2798.182 },
    @{ Cell = "K77"; Value = 30059.6875 },
    @{ Cell = "L77"; Value = 13990.91 },
    @{ Cell = "M77"; Value = -25691.6875 },
    @{ Cell = "N77"; Value = -22726.91 },
    @{ Cell = "H136"; Value = 2104.6843 },
    @{ Cell = "I136"; Value = 1326.7778 },
    @{ Cell = "J136"; Value = 2804.8 },
    @{ Cell = "K136"; Value = 3980.3334 },
    @{ Cell = "L136"; Value = 8414.400000000001 },
    @{ Cell = "M136"; Value = -1430.3334 },
    @{ Cell = "N136"; Value = -13514.4 },
    @{ Cell = "H137"; Value = 40135 },
    @{ Cell = "J137"; Value = 40135 },
    @{ Cell = "L137"; Value = 40135 },
    @{ Cell = "N137"; Value = -50335 }
)

$sheetEdits["BSM"] = @(
    @{ Cell = "H86"; Value = 1585.3572 },
    @{ Cell = "I86"; Value = 1576.4445 },
    @{ Cell = "J86"; Value = 1601.4 },
    @{ Cell = "K86"; Value = 1576.4445 },
    @{ Cell = "L86"; Value = 1601.4 },
    @{ Cell = "M86"; Value = -453.4445000000001 },
    @{ Cell = "N86"; Value = -3847.4 },
    @{ Cell = "H89"; Value = 1585.3572 },
    @{ Cell = "I89"; Value = 1576.4445 },
    @{ Cell = "J89"; Value = 1601.4 },
    @{ Cell = "K89"; Value = 7882.2225 },
    @{ Cell = "L89"; Value = 8007 },
    @{ Cell = "M89"; Value = -2266.2225 },
    @{ Cell = "N89"; Value = -19239 },
    @{ Cell = "H94"; Value = 815.7083 },
    @{ Cell = "I94"; Value = 652.4 },
    @{ Cell = "J94"; Value = 1632.25 },
    @{ Cell = "K94"; Value = 652.4 },
    @{ Cell = "L94"; Value = 1632.25 },
    @{ Cell = "M94"; Value = -201.4 },
    @{ Cell = "N94"; Value = -2534.25 },
    @{ Cell = "H99"; Value = 3506.25 },
    @{ Cell = "I99"; Value = 1100 },
    @{ Cell = "J99"; Value = 5377.778 },
    @{ Cell = "K99"; Value = 1100 },
    @{ Cell = "L99"; Value = 5377.778 },
    @{ Cell = "M99"; Value = 398 },
    @{ Cell = "N99"; Value = -8373.778 },
    @{ Cell = "H137"; Value = 43486.668 },
    @{ Cell = "J137"; Value = 43486.668 },
    @{ Cell = "L137"; Value = 43486.668 },
    @{ Cell = "N137"; Value = -53686.668 }
)

$sheetEdits["CRP"] = @(
    @{ Cell = "H31"; Value = 2699.5134 },
    @{ Cell = "I31"; Value = 1017.7222 },
    @{ Cell = "J31"; Value = 4292.7896 },
    @{ Cell = "K31"; Value = 1017.7222 },
    @{ Cell = "L31"; Value = 4292.7896 },
    @{ Cell = "M31"; Value = -722.7222 },
    @{ Cell = "N31"; Value = -4882.7896 },
    @{ Cell = "H34"; Value = 2699.5134 },
    @{ Cell = "I34"; Value = 1017.7222 },
    @{ Cell = "J34"; Value = 4292.7896 },
    @{ Cell = "K34"; Value = 1017.7222 },
    @{ Cell = "L34"; Value = 4292.7896 },
    @{ Cell = "M34"; Value = -815.7222 },
    @{ Cell = "N34"; Value = -4696.7896 },
    @{ Cell = "H60"; Value = 19154.424 },
    @{ Cell = "J60"; Value = 19154.424 },
    @{ Cell = "L60"; Value = 19154.424 },
    @{ Cell = "N60"; Value = -20176.424 },
    @{ Cell = "H105"; Value = 1848.6842 },
    @{ Cell = "I105"; Value = 1475.6666 },
    @{ Cell = "J105"; Value = 3247.5 },
    @{ Cell = "K105"; Value = 1475.6666 },
    @{ Cell = "L105"; Value = 3247.5 },
    @{ Cell = "M105"; Value = 271.3334 },
    @{ Cell = "N105"; Value = -6741.5 },
    @{ Cell = "H110"; Value = 49990 },
    @{ Cell = "J110"; Value = 49990 },
    @{ Cell = "L110"; Value = 49990 },
    @{ Cell = "N110"; Value = -58170 }
)

$sheetEdits["CUL"] = @(
    @{ Cell = "H5"; Value = 557318.3 },
    @{ Cell = "J5"; Value = 835705.2 },
    @{ Cell = "L5"; Value = 2507115.6 },
    @{ Cell = "N5"; Value = -2507339.6 },
    @{ Cell = "H80"; Value = 18479.1 },
    @{ Cell = "J80"; Value = 18479.1 },
    @{ Cell = "L80"; Value = 55437.3 },
    @{ Cell = "N80"; Value = -57309.3 },
    @{ Cell = "H83"; Value = 18479.1 },
    @{ Cell = "J83"; Value = 18479.1 },
    @{ Cell = "L83"; Value = 166311.9 },
    @{ Cell = "N83"; Value = -175671.9 },
    @{ Cell = "H113"; Value = 4032860.8 },
    @{ Cell = "I113"; Value = 619.64703 },
    @{ Cell = "J113"; Value = 8929153 },
    @{ Cell = "K113"; Value = 1858.94109 },
    @{ Cell = "L113"; Value = 26787459 },
    @{ Cell = "M113"; Value = 311.0589100000002 },
    @{ Cell = "N113"; Value = -26791799 },
    @{ Cell = "H129"; Value = 3668.0908 },
    @{ Cell = "J129"; Value = 3183.1667 },
    @{ Cell = "L129"; Value = 9549.500100000001 },
    @{ Cell = "N129"; Value = -19549.5001 },
    @{ Cell = "H132"; Value = 2406.95 },
    @{ Cell = "I132"; Value = 966.6667 },
    @{ Cell = "J132"; Value = 3585.3635 },
    @{ Cell = "K132"; Value = 8700.0003 },
    @{ Cell = "L132"; Value = 32268.2715 },
    @{ Cell = "M132"; Value = -6170.0003 },
    @{ Cell = "N132"; Value = -37328.2715 },
    @{ Cell = "H135"; Value = 557318.3 },
    @{ Cell = "J135"; Value = 835705.2 },
    @{ Cell = "L135"; Value = 7521346.8 },
    @{ Cell = "N135"; Value = -7526416.8 }
)

$sheetEdits["GSM"] = @(
    @{ Cell = "H80"; Value = 35717056 },
    @{ Cell = "I80"; Value = 62502676 },
    @{ Cell = "J80"; Value = 2900 },
    @{ Cell = "K80"; Value = 62502676 },
    @{ Cell = "L80"; Value = 2900 },
    @{ Cell = "M80"; Value = -62501678 },
    @{ Cell = "N80"; Value = -4896 },
    @{ Cell = "H83"; Value = 35717056 },
    @{ Cell = "I83"; Value = 62502676 },
    @{ Cell = "J83"; Value = 2900 },
    @{ Cell = "K83"; Value = 312513380 },
    @{ Cell = "L83"; Value = 14500 },
    @{ Cell = "M83"; Value = -312508388 },
    @{ Cell = "N83"; Value = -24484 },
    @{ Cell = "H97"; Value = 1042.3334 },
    @{ Cell = "I97"; Value = 1064.25 },
    @{ Cell = "J97"; Value = 998.5 },
    @{ Cell = "K97"; Value = 1064.25 },
    @{ Cell = "L97"; Value = 998.5 },
    @{ Cell = "M97"; Value = -568.25 },
    @{ Cell = "N97"; Value = -1990.5 },
    @{ Cell = "H132"; Value = 3307.72 },
    @{ Cell = "I132"; Value = 1893 },
    @{ Cell = "J132"; Value = 5429.8 },
    @{ Cell = "K132"; Value = 5679 },
    @{ Cell = "L132"; Value = 16289.4 },
    @{ Cell = "M132"; Value = -3149 },
    @{ Cell = "N132"; Value = -21349.4 },
    @{ Cell = "H137"; Value = 37212 },
    @{ Cell = "J137"; Value = 37212 },
    @{ Cell = "L137"; Value = 37212 },
    @{ Cell = "N137"; Value = -47412 }
)

$sheetEdits["LTW"] = @(
    @{ Cell = "H40"; Value = 7988 },
    @{ Cell = "I40"; Value = 6484 },
    @{ Cell = "J40"; Value = 12500 },
    @{ Cell = "K40"; Value = 6484 },
    @{ Cell = "L40"; Value = 12500 },
    @{ Cell = "M40"; Value = -6348 },
    @{ Cell = "N40"; Value = -12772 },
    @{ Cell = "H46"; Value = 2965.5293 },
    @{ Cell = "I46"; Value = 3050 },
    @{ Cell = "J46"; Value = 2939.5386 },
    @{ Cell = "K46"; Value = 3050 },
    @{ Cell = "L46"; Value = 2939.5386 },
    @{ Cell = "M46"; Value = -2862 },
    @{ Cell = "N46"; Value = -3315.5386 },
    @{ Cell = "H68"; Value = 656.32 },
    @{ Cell = "I68"; Value = 656.32 },
    @{ Cell = "K68"; Value = 656.32 },
    @{ Cell = "M68"; Value = 92.67999999999995 },
    @{ Cell = "H71"; Value = 656.32 },
    @{ Cell = "I71"; Value = 656.32 },
    @{ Cell = "K71"; Value = 3281.6 },
    @{ Cell = "M71"; Value = 462.3999999999996 }
)

$sheetEdits["WVR"] = @(
    @{ Cell = "H57"; Value = 17725 },
    @{ Cell = "J57"; Value = 17725 },
    @{ Cell = "L57"; Value = 17725 },
    @{ Cell = "N57"; Value = -19233 },
    @{ Cell = "H96"; Value = 142930700 },
    @{ Cell = "J96"; Value = 5495 },
    @{ Cell = "L96"; Value = 5495 },
    @{ Cell = "N96"; Value = -8241 },
    @{ Cell = "H122"; Value = 7114.143 },
    @{ Cell = "I122"; Value = 3699.75 },
    @{ Cell = "K122"; Value = 11099.25 },
    @{ Cell = "M122"; Value = -8649.25 }
)

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($edit in $sheetEdits[$sheetName]) {
        $ws.Range($edit.Cell).Value = $edit.Value
    }
}

